$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, B, C, D) for rows 344-357,
# corresponding to dates 2021-08-10 through 2021-08-23
# (update al 23 agosto 2021).
$newRows = @(
    @(44418, 0, 1, 145.7725947521866),
    @(44419, 0, 1, 145.7725947521866),
    @(44420, 0, 1, 145.7725947521866),
    @(44421, 0, 1, 145.7725947521866),
    @(44422, 0, 1, 145.7725947521866),
    @(44423, 0, 0, 0),
    @(44424, 1, 1, 145.7725947521866),
    @(44425, 0, 1, 145.7725947521866),
    @(44426, 0, 1, 145.7725947521866),
    @(44427, 1, 2, 291.5451895043732),
    @(44428, 0, 2, 291.5451895043732),
    @(44429, 0, 2, 291.5451895043732),
    @(44430, 0, 2, 291.5451895043732),
    @(44431, 0, 1, 145.7725947521866)
)

$lastExistingRow = 343
$startRow = $lastExistingRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Copy formatting (style/number format/borders) from the row above
    # so the new row matches the look of the preceding data rows.
    $srcRange = $ws.Range("A" + ($r - 1) + ":D" + ($r - 1))
    $dstRange = $ws.Range("A" + $r + ":D" + $r)
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($r, 1).Value2 = $data[0]
    $ws.Cells.Item($r, 2).Value2 = $data[1]
    $ws.Cells.Item($r, 3).Value2 = $data[2]
    $ws.Cells.Item($r, 4).Value2 = $data[3]
}
